$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.270.11"
$ws.Range("E2").Value = "  -2.09%  "
$ws.Range("D3").Value = "3.489.86"
$ws.Range("E3").Value = "  -2.04%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "'612.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.08%  "
$ws.Range("D6").Value = "'185.62"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.23%  "
$ws.Range("D7").Value = "'0.631"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "'0.217"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.17%  "
$ws.Range("D10").Value = "'0.651"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.28%  "
$ws.Range("D11").Value = "'53.07"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.46%  "
$ws.Range("D12").Value = "'0.0000308"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.71%  "
$ws.Range("D13").Value = "'9.56"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.69%  "
$ws.Range("D14").Value = "4.044.54"
$ws.Range("E14").Value = "  -1.94%  "
$ws.Range("D15").Value = "'602.61"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.83%  "
$ws.Range("D16").Value = "69.314.18"
$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D17").Value = "'18.89"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.00%  "
$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D18").Value = "'12.58"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.48%  "
$ws.Range("D19").Value = "3.484.91"
$ws.Range("E19").Value = "  -1.99%  "
$ws.Range("E20").Value = "  -0.28%  "
$ws.Range("D21").Value = "'0.989"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.55%  "
$ws.Range("D22").Value = "'17.24"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.39%  "
$ws.Range("D23").Value = "'104.79"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +9.98%  "
$ws.Range("D24").Value = "'4.66"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.41%  "
$ws.Range("D25").Value = "'5.07"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.19%  "
$ws.Range("D26").Value = "'3.03"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.02%  "
$ws.Range("D27").Value = "'11.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.72%  "
$ws.Range("D28").Value = "'10.02"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +9.73%  "
$ws.Range("D29").Value = "'33.49"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.13%  "
$ws.Range("D30").Value = "'6.97"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.57%  "
$ws.Range("D31").Value = "'12.45"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.24%  "
$ws.Range("E32").Value = "  -0.28%  "
$ws.Range("D33").Value = "'3.86"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +15.03%  "
$ws.Range("D34").Value = "'63.34"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.02%  "
$ws.Range("D35").Value = "'3.18"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -7.67%  "
$ws.Range("D36").Value = "'0.999"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.04%  "
$ws.Range("D37").Value = "'520.31"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.50%  "
$ws.Range("D38").Value = "'0.399"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.03%  "
$ws.Range("D39").Value = "3.596.56"
$ws.Range("E39").Value = "  +0.80%  "
$ws.Range("D40").Value = "'3.60"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.39%  "
$ws.Range("D41").Value = "'36.79"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.78%  "
$ws.Range("D42").Value = "0.0₃0774"
$ws.Range("E42").Value = "  -2.42%  "
$ws.Range("E43").Value = "  +1.12%  "
$ws.Range("D44").Value = "'0.0462"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.50%  "
$ws.Range("D45").Value = "'2.96"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.07%  "
$ws.Range("E46").Value = "  +3.61%  "
$ws.Range("D47").Value = "'3.32"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.07%  "
$ws.Range("D48").Value = "'8.83"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.69%  "
$ws.Range("E49").Value = "  +0.32%  "
$ws.Range("D50").Value = "'0.000244"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.81%  "
$ws.Range("B51").Value = "OceanProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/aAKLSV5-0+oceanprotocol-ocean"
$ws.Range("D51").Value = "'1.37"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -8.97%  "
